$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-04 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-05 Friday", 2) | Out-Null
$d.Content.Find.Execute("47×90=4230", $true, $false, $false, $false, $false, $true, 1, $false, "20×87=1740", 2) | Out-Null
$d.Content.Find.Execute("57×61=3477", $true, $false, $false, $false, $false, $true, 1, $false, "27×48=1296", 2) | Out-Null
$d.Content.Find.Execute("59×67=3953", $true, $false, $false, $false, $false, $true, 1, $false, "75×41=3075", 2) | Out-Null
$d.Content.Find.Execute("78×83=6474", $true, $false, $false, $false, $false, $true, 1, $false, "31×78=2418", 2) | Out-Null
$d.Content.Find.Execute("40×73=2920", $true, $false, $false, $false, $false, $true, 1, $false, "91×49=4459", 2) | Out-Null
$d.Content.Find.Execute("53×76=4028", $true, $false, $false, $false, $false, $true, 1, $false, "95×97=9215", 2) | Out-Null
$d.Content.Find.Execute("82×97=7954", $true, $false, $false, $false, $false, $true, 1, $false, "42×45=1890", 2) | Out-Null
$d.Content.Find.Execute("40×13=520", $true, $false, $false, $false, $false, $true, 1, $false, "76×31=2356", 2) | Out-Null
$d.Content.Find.Execute("54×77=4158", $true, $false, $false, $false, $false, $true, 1, $false, "31×36=1116", 2) | Out-Null
$d.Content.Find.Execute("91×42=3822", $true, $false, $false, $false, $false, $true, 1, $false, "37×32=1184", 2) | Out-Null
$d.Content.Find.Execute("39×52=2028", $true, $false, $false, $false, $false, $true, 1, $false, "27×69=1863", 2) | Out-Null
$d.Content.Find.Execute("50×15=750", $true, $false, $false, $false, $false, $true, 1, $false, "76×78=5928", 2) | Out-Null
$d.Content.Find.Execute("49×99=4851", $true, $false, $false, $false, $false, $true, 1, $false, "34×90=3060", 2) | Out-Null
$d.Content.Find.Execute("36×24=864", $true, $false, $false, $false, $false, $true, 1, $false, "12×47=564", 2) | Out-Null
$d.Content.Find.Execute("46×99=4554", $true, $false, $false, $false, $false, $true, 1, $false, "82×81=6642", 2) | Out-Null
$d.Content.Find.Execute("50×78=3900", $true, $false, $false, $false, $false, $true, 1, $false, "72×20=1440", 2) | Out-Null
$d.Content.Find.Execute("98×47=4606", $true, $false, $false, $false, $false, $true, 1, $false, "76×24=1824", 2) | Out-Null
$d.Content.Find.Execute("67×44=2948", $true, $false, $false, $false, $false, $true, 1, $false, "71×53=3763", 2) | Out-Null
$d.Content.Find.Execute("81×99=8019", $true, $false, $false, $false, $false, $true, 1, $false, "80×75=6000", 2) | Out-Null
$d.Content.Find.Execute("41×17=697", $true, $false, $false, $false, $false, $true, 1, $false, "82×91=7462", 2) | Out-Null
$d.Content.Find.Execute("29×48=1392", $true, $false, $false, $false, $false, $true, 1, $false, "41×14=574", 2) | Out-Null
$d.Content.Find.Execute("12×13=156", $true, $false, $false, $false, $false, $true, 1, $false, "24×23=552", 2) | Out-Null
$d.Content.Find.Execute("85×14=1190", $true, $false, $false, $false, $false, $true, 1, $false, "19×67=1273", 2) | Out-Null
$d.Content.Find.Execute("81×89=7209", $true, $false, $false, $false, $false, $true, 1, $false, "79×76=6004", 2) | Out-Null
